$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small single-cell value corrections (Q column re-computations, SIMBAD update) ---
$ws.Cells.Item(9, 17).Value  = 1.616889990960015
$ws.Cells.Item(29, 17).Value = 2.608985367019665
$ws.Cells.Item(79, 17).Value = 5.029247845400172
$ws.Cells.Item(80, 17).Value = 5.06680867658183
$ws.Cells.Item(92, 17).Value = 33.7642303694855
$ws.Cells.Item(93, 17).Value = 187.1954334784243

# --- Rows 32-34 were reordered (2443.01 moved after the two 6965.xx rows),
#     and the Q value for TOI 2443.01 was recomputed. ---

# New row 32 (previously row 33: TOI 6965.01) - unchanged values, just moved up
$ws.Cells.Item(32, 1).Value  = 6965.01
$ws.Cells.Item(32, 2).Value  = 80224448
$ws.Cells.Item(32, 3).Value  = 103.678273
$ws.Cells.Item(32, 4).Value  = 24.245141
$ws.Cells.Item(32, 5).Value  = 2459505.193695
$ws.Cells.Item(32, 6).Value  = 5.9693397
$ws.Cells.Item(32, 7).Value  = 3.441038
$ws.Cells.Item(32, 8).Value  = 105.5424205
$ws.Cells.Item(32, 9).Value  = 1.3473608
$ws.Cells.Item(32, 10).Value = 334.6170929
$ws.Cells.Item(32, 11).Value = 1090.8292631
$ws.Cells.Item(32, 12).Value = 6.3103
$ws.Cells.Item(32, 13).Value = 31.1664
$ws.Cells.Item(32, 14).Value = 6007
$ws.Cells.Item(32, 15).Value = 4.37729
$ws.Cells.Item(32, 16).Value = 1.12993
$ws.Cells.Item(32, 17).Value = 2.880011854113967
$ws.Cells.Item(32, 18).Value = 0.05173693485062012

# New row 33 (previously row 34: TOI 6965.02) - unchanged values, just moved up
$ws.Cells.Item(33, 1).Value  = 6965.02
$ws.Cells.Item(33, 2).Value  = 80224448
$ws.Cells.Item(33, 3).Value  = 103.678273
$ws.Cells.Item(33, 4).Value  = 24.245141
$ws.Cells.Item(33, 5).Value  = 2459508.008646
$ws.Cells.Item(33, 6).Value  = 28.0693949
$ws.Cells.Item(33, 7).Value  = 6.0393818
$ws.Cells.Item(33, 8).Value  = 154.3008171
$ws.Cells.Item(33, 9).Value  = 1.6221344
$ws.Cells.Item(33, 10).Value = 42.4755605
$ws.Cells.Item(33, 11).Value = 651.1102456
$ws.Cells.Item(33, 12).Value = 6.3103
$ws.Cells.Item(33, 13).Value = 31.1664
$ws.Cells.Item(33, 14).Value = 6007
$ws.Cells.Item(33, 15).Value = 4.37729
$ws.Cells.Item(33, 16).Value = 1.12993
$ws.Cells.Item(33, 17).Value = 2.880011854113967
$ws.Cells.Item(33, 18).Value = 0.05173693485062012

# New row 34 (previously row 32: TOI 2443.01) - moved down, Q recomputed (2.839885605748058 -> 2.970022591198855)
$ws.Cells.Item(34, 1).Value  = 2443.01
$ws.Cells.Item(34, 2).Value  = 318753380
$ws.Cells.Item(34, 3).Value  = 40.179861
$ws.Cells.Item(34, 4).Value  = 1.199676
$ws.Cells.Item(34, 5).Value  = 2459148.098617
$ws.Cells.Item(34, 6).Value  = 15.6692322
$ws.Cells.Item(34, 7).Value  = 4.562853
$ws.Cells.Item(34, 8).Value  = 1393.3814802
$ws.Cells.Item(34, 9).Value  = 2.6868893
$ws.Cells.Item(34, 10).Value = 13.2931176
$ws.Cells.Item(34, 11).Value = 486.9971964
$ws.Cells.Item(34, 12).Value = 8.296900000000001
$ws.Cells.Item(34, 13).Value = 23.9258
$ws.Cells.Item(34, 14).Value = 4214.44
$ws.Cells.Item(34, 15).Value = 4.52845
$ws.Cells.Item(34, 16).Value = 0.732115
$ws.Cells.Item(34, 17).Value = 2.970022591198855
$ws.Cells.Item(34, 18).Value = 0.1818675256571092

$wb.Save()
